# Updates the "Estado de Cuenta" worksheet: refresh the quota/mora summary
# figures, replace the detail table with the new worker records (part 1 of
# the new statement), and drop the now-obsolete trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 11/13) -------------------------------------------
$ws.Range("E11").Value = 268063
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 3

# --- Detail table (rows 16-25) ---------------------------------------------
# Column layout: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
# E=Periodo Mora, F=Valor Mora, G=Salario Basico
$data = @(
    @("CC", "33337424",   "MARBEL LUZ BALLESTAS BUELVAS",     "2010", 35112, 908526),
    @("CC", "1049942967", "ANAYIBIS PEREZ HERNANDEZ",         "2010", 35112, 877803),
    @("CC", "1049939325", "MARIA VANESSA BENAVIDES MARIMON",  "2010", 35112, 877803),
    @("CC", "41371858",   "MARIA TERESA PADILLA CAÃ?ATE",     "2010", 35112, 908526),
    @("CC", "33273074",   "KATIA MERCEDES CONTRERAS ARDILA",  "2011", 35112, 908526),
    @("CC", "33337424",   "MARBEL LUZ BALLESTAS BUELVAS",     "2011", 35112, 908526),
    @("CC", "1049942967", "ANAYIBIS PEREZ HERNANDEZ",         "2011", 10534, 877803),
    @("CC", "1049939325", "MARIA VANESSA BENAVIDES MARIMON",  "2011", 10534, 877803),
    @("CC", "41371858",   "MARIA TERESA PADILLA CAÃ?ATE",     "2011", 35112, 908526),
    @("CC", "41371858",   "MARIA TERESA PADILLA CAÃ?ATE",     "2102",  1211, 908526)
)

$row = 16
foreach ($rec in $data) {
    $ws.Range("B$row").Value = $rec[0]
    $ws.Range("C$row").Value = $rec[1]
    $ws.Range("D$row").Value = $rec[2]
    $ws.Range("E$row").Value = $rec[3]
    $ws.Range("F$row").Value = $rec[4]
    $ws.Range("G$row").Value = $rec[5]
    $row++
}

# --- Drop the four rows that no longer apply (old rows 26-29) --------------
# This also shifts the trailing signature block (old rows 34/35) up to 30/31.
$ws.Range("26:29").EntireRow.Delete()
